$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 190; existing rows 190-218 shift down to 191-219.
$ws.Rows(190).Insert()

# Populate the newly inserted row 190 with the new weekly record.
$ws.Cells.Item(190, 1).Value = 5
$ws.Cells.Item(190, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(190, 3).Value = "Maule"
$ws.Cells.Item(190, 4).Value = "2021-11-05"
$ws.Cells.Item(190, 5).Value = 7
$ws.Cells.Item(190, 6).Value = 100112023
$ws.Cells.Item(190, 7).Value = "Brócoli"
$ws.Cells.Item(190, 8).Value = "Sin especificar"
$ws.Cells.Item(190, 9).Value = "Primera"
$ws.Cells.Item(190, 10).Value = 6000
$ws.Cells.Item(190, 11).Value = 500
$ws.Cells.Item(190, 12).Value = 500
$ws.Cells.Item(190, 13).Value = 500
$ws.Cells.Item(190, 14).Value = "`$/unidad"
$ws.Cells.Item(190, 15).Value = "Región del Maule"
$ws.Cells.Item(190, 16).Value = 500
$ws.Cells.Item(190, 17).Value = 1
$ws.Cells.Item(190, 18).Value = "Hortaliza"
